$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.568.31"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.962.42"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  -0.06%  "
$r = $ws.Range("D6")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "0.617"
$r.Style = $sty
$ws.Range("E6").Value = "  -0.60%  "
$r = $ws.Range("D7")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "58.67"
$r.Style = $sty
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E9").Value = "  +2.97%  "
$r = $ws.Range("D10")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "0.0806"
$r.Style = $sty
$ws.Range("E10").Value = "  -6.77%  "
$ws.Range("E11").Value = "  -0.87%  "
$r = $ws.Range("D12")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "22.09"
$r.Style = $sty
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "2.249.89"
$ws.Range("E14").Value = "  +0.19%  "
$r = $ws.Range("D15")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "13.79"
$r.Style = $sty
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").Value = "1.964.12"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "36.543.61"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").Value = "0.0₃0858"
$ws.Range("E20").Value = "  -2.57%  "
$r = $ws.Range("D21")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "228.74"
$r.Style = $sty
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("E23").Value = "  -0.03%  "
$r = $ws.Range("D24")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "2.44"
$r.Style = $sty
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("E25").Value = "  +1.66%  "
$r = $ws.Range("D26")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "9.26"
$r.Style = $sty
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("E27").Value = "  +1.31%  "
$r = $ws.Range("D28")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "160.67"
$r.Style = $sty
$ws.Range("E28").Value = "  -0.95%  "
$r = $ws.Range("D29")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "19.43"
$r.Style = $sty
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("E31").Value = "  -3.02%  "
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("E33").Value = "  -3.95%  "
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("E37").Value = "  +11.37%  "
$ws.Range("E38").Value = "  -0.18%  "
$r = $ws.Range("D39")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "5.78"
$r.Style = $sty
$ws.Range("E39").Value = "  -10.01%  "
$r = $ws.Range("D40")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "0.0983"
$r.Style = $sty
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("E42").Value = "  -0.80%  "
$r = $ws.Range("D43")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "0.0211"
$r.Style = $sty
$ws.Range("E43").Value = "  -0.19%  "
$r = $ws.Range("D44")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "16.00"
$r.Style = $sty
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").Value = "1.365.03"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "2.141.28"
$ws.Range("E50").Value = "  +0.23%  "
$r = $ws.Range("D51")
$sty = $r.Style
$r.NumberFormat = "@"
$r.Value = "43.79"
$r.Style = $sty
$ws.Range("E51").Value = "  -5.17%  "
